$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing "Ammount" quantity for TB6612FNG (row 2) from 2 to 5
$ws.Range("D2").Value = 5

# Add a new component row (row 6) for the LM1086 voltage regulator
$ws.Range("A6").Value = "LM1086"
$ws.Range("B6").Value = "TO-220"
$ws.Range("C6").Value = "LM1086 Stab. 5v - Mouser"
$ws.Range("D6").Value = 1

# Hyperlink the new datasheet/shop link, then restore the standard
# "hyperlink" cell look used by the other rows in column C
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.mouser.com/ProductDetail/Texas-Instruments/LM1086CT-5.0-NOPB")
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the author last left it
[void]$ws.Range("D8").Select()
